$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 126
$ws.Range("I3").Value = 196
$ws.Range("J3").Value = 240
$ws.Range("C6").Value = 491
$ws.Range("D6").Value = 426
$ws.Range("E6").Value = 490
$ws.Range("F6").Value = 557
$ws.Range("G6").Value = 440
$ws.Range("H6").Value = 453
$ws.Range("I6").Value = 509
$ws.Range("C7").Value = 648
$ws.Range("D7").Value = 666
$ws.Range("E7").Value = 724
$ws.Range("F7").Value = 805
$ws.Range("G7").Value = 674
$ws.Range("H7").Value = 737
$ws.Range("I7").Value = 849
$ws.Range("J7").Value = 814

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("F6").Value = 52
$ws.Range("H6").Value = 38
$ws.Range("F7").Value = 61
$ws.Range("H7").Value = 49

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("D5").Value = 11
$ws.Range("D6").Value = 14

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 10

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 9
$ws.Range("F6").Value = 39
$ws.Range("F7").Value = 59
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G10").Value = 5
$ws.Range("I19").Value = 23
$ws.Range("F28").Value = 59
$ws.Range("J28").Value = 35
$ws.Range("F32").Value = 61
$ws.Range("H32").Value = 49
$ws.Range("I47").Value = 26
$ws.Range("D50").Value = 14
$ws.Range("H53").Value = 107
$ws.Range("J53").Value = 126
$ws.Range("D65").Value = 28
$ws.Range("F65").Value = 40
$ws.Range("F77").Value = 23
$ws.Range("E88").Value = 10
$ws.Range("C94").Value = 6
$ws.Range("C98").Value = 648
$ws.Range("D98").Value = 666
$ws.Range("E98").Value = 724
$ws.Range("F98").Value = 805
$ws.Range("G98").Value = 674
$ws.Range("H98").Value = 737
$ws.Range("I98").Value = 849
$ws.Range("J98").Value = 814

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 39
$ws.Range("H6").Value = 68
$ws.Range("H7").Value = 107
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("D5").Value = 27
$ws.Range("F5").Value = 33
$ws.Range("D6").Value = 28
$ws.Range("F6").Value = 40

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 23

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 5

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I3").Value = 2
$ws.Range("I6").Value = 26

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("C5").Value = 5
$ws.Range("C6").Value = 6
